$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.318.24"
$ws.Range("E2").Value = "  -0.08%  "

$ws.Range("D3").Value = "3.373.27"
$ws.Range("E3").Value = "  +1.60%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "'573.82"
$ws.Range("E5").Value = "  +1.18%  "

$ws.Range("D6").Value = "'137.18"
$ws.Range("E6").Value = "  +6.50%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").Value = "3.373.56"
$ws.Range("E8").Value = "  +1.58%  "

$ws.Range("D9").Value = "'0.477"
$ws.Range("E9").Value = "  -0.71%  "

$ws.Range("D10").Value = "'7.62"
$ws.Range("E10").Value = "  +4.24%  "

$ws.Range("E11").Value = "  +4.31%  "

$ws.Range("E12").Value = "  +3.71%  "

$ws.Range("D13").Value = "3.940.10"
$ws.Range("E13").Value = "  +1.40%  "

$ws.Range("E14").Value = "  +1.86%  "

$ws.Range("D15").Value = "'0.0000177"
$ws.Range("E15").Value = "  +4.40%  "

$ws.Range("D16").Value = "3.365.53"
$ws.Range("E16").Value = "  +1.55%  "

$ws.Range("D17").Value = "'25.33"
$ws.Range("E17").Value = "  +2.34%  "

$ws.Range("D18").Value = "61.305.45"
$ws.Range("E18").Value = "  -0.28%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'13.93"
$ws.Range("E19").Value = "  +3.72%  "

$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'5.88"
$ws.Range("E20").Value = "  +3.23%  "

$ws.Range("D21").Value = "'9.36"
$ws.Range("E21").Value = "  +3.73%  "

$ws.Range("D22").Value = "'380.51"
$ws.Range("E22").Value = "  +7.00%  "

$ws.Range("D23").Value = "'0.568"
$ws.Range("E23").Value = "  +1.67%  "

$ws.Range("D24").Value = "3.501.92"
$ws.Range("E24").Value = "  +1.44%  "

$ws.Range("E25").Value = "  -0.09%  "

$ws.Range("D26").Value = "'70.56"
$ws.Range("E26").Value = "  +1.49%  "

$ws.Range("D27").Value = "'0.0000122"
$ws.Range("E27").Value = "  +12.51%  "

$ws.Range("D28").Value = "'1.62"
$ws.Range("E28").Value = "  +11.23%  "

$ws.Range("D29").Value = "'7.80"
$ws.Range("E29").Value = "  +7.27%  "

$ws.Range("D30").Value = "'0.995"
$ws.Range("E30").Value = "  -0.58%  "

$ws.Range("D31").Value = "'8.25"
$ws.Range("E31").Value = "  +4.62%  "

$ws.Range("D32").Value = "'0.157"
$ws.Range("E32").Value = "  +4.69%  "

$ws.Range("E33").Value = "  +0.63%  "

$ws.Range("E34").Value = "  -0.07%  "

$ws.Range("D35").Value = "3.395.11"
$ws.Range("E35").Value = "  +1.41%  "

$ws.Range("D36").Value = "'23.46"
$ws.Range("E36").Value = "  +3.63%  "

$ws.Range("D37").Value = "'5.50"
$ws.Range("E37").Value = "  +2.71%  "

$ws.Range("D38").Value = "'7.05"
$ws.Range("E38").Value = "  +2.86%  "

$ws.Range("E39").Value = "  +3.32%  "

$ws.Range("D40").Value = "'160.85"
$ws.Range("E40").Value = "  -0.27%  "

$ws.Range("D41").Value = "'0.0789"
$ws.Range("E41").Value = "  +3.36%  "

$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.12%  "

$ws.Range("D43").Value = "'1.73"
$ws.Range("E43").Value = "  +11.07%  "

$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D44").Value = "'1.22"
$ws.Range("E44").Value = "  +8.47%  "

$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "'4.43"
$ws.Range("E45").Value = "  +1.31%  "

$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.765"
$ws.Range("E46").Value = "  +2.76%  "

$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "'41.44"
$ws.Range("E47").Value = "  +0.95%  "

$ws.Range("D48").Value = "'23.26"
$ws.Range("E48").Value = "  +4.38%  "

$ws.Range("E49").Value = "  +3.41%  "

$ws.Range("D50").Value = "'22.78"
$ws.Range("E50").Value = "  +6.29%  "

$ws.Range("D51").Value = "2.335.80"
$ws.Range("E51").Value = "  +7.96%  "
